$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$json = @'
{
    "$schema": "http:\/\/json-schema.org\/draft-07\/schema#",
    "description": "MOSIP Sample identity",
    "additionalProperties": false,
    "title": "MOSIP identity",
    "type": "object",
    "definitions": {
        "simpleType": {
            "uniqueItems": true,
            "additionalItems": false,
            "type": "array",
            "items": {
                "additionalProperties": false,
                "type": "object",
                "required": [
                    "language",
                    "value"
                ],
                "properties": {
                    "language": {
                        "type": "string"
                    },
                    "value": {
                        "type": "string"
                    }
                }
            }
        },
        "taggedListType": {
            "uniqueItems": true,
            "additionalItems": false,
            "type": "array",
            "items": {
                "additionalProperties": false,
                "type": "object",
                "required": [
                    "value"
                ],
                "properties": {
                    "value": {
                        "type": "string"
                    },
                    "tags": {
                        "uniqueItems": true,
                        "type": "array",
                        "items": {
                            "type": "string",
                            "enum": [
                                "notification",
                                "handle"
                            ]
                        }
                    }
                }
            }
        },
        "hashType": {
            "additionalProperties": false,
            "type": "object",
            "properties": {
                "hash": {
                    "type": "string"
                },
                "salt": {
                    "type": "string"
                }
            }
        },
        "documentType": {
            "additionalProperties": false,
            "type": "object",
            "properties": {
                "format": {
                    "type": "string"
                },
                "type": {
                    "type": "string"
                },
                "value": {
                    "type": "string"
                },
                "refNumber": {
                    "type": [
                        "string",
                        "null"
                    ]
                }
            }
        },
        "biometricsType": {
            "additionalProperties": false,
            "type": "object",
            "properties": {
                "format": {
                    "type": "string"
                },
                "version": {
                    "type": "number",
                    "minimum": 0
                },
                "value": {
                    "type": "string"
                }
            }
        }
    },
    "properties": {
        "identity": {
            "additionalProperties": false,
            "type": "object",
            "required": [
                "IDSchemaVersion",
                "fullName",
                "dateOfBirth",
                "gender",
                "addressLine1",
                "addressLine2",
                "addressLine3",
                "region",
                "province",
                "city",
                "zone",
                "postalCode",
                "phone",
                "email",
                "proofOfIdentity",
                "individualBiometrics"
            ],
            "properties": {
                "proofOfAddress": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/documentType"
                },
                "gender": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "city": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{0,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "postalCode": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^[(?i)A-Z0-9]{5}$|^NA$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "proofOfException-1": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/documentType"
                },
                "referenceIdentityNumber": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^([0-9]{10,30})$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "kyc",
                    "type": "string",
                    "fieldType": "default"
                },
                "individualBiometrics": {
                    "bioAttributes": [
                        "leftEye",
                        "rightEye",
                        "rightIndex",
                        "rightLittle",
                        "rightRing",
                        "rightMiddle",
                        "leftIndex",
                        "leftLittle",
                        "leftRing",
                        "leftMiddle",
                        "leftThumb",
                        "rightThumb",
                        "face"
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/biometricsType"
                },
                "province": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{0,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "zone": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "proofOfDateOfBirth": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/documentType"
                },
                "addressLine1": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{0,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "addressLine2": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{3,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "residenceStatus": {
                    "bioAttributes": [],
                    "fieldCategory": "kyc",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "addressLine3": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{3,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "email": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^[A-Za-z0-9_\\-]+(\\.[A-Za-z0-9_]+)*@[A-Za-z0-9_-]+(\\.[A-Za-z0-9_]+)*(\\.[a-zA-Z]{2,})$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "introducerRID": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "introducerBiometrics": {
                    "bioAttributes": [
                        "leftEye",
                        "rightEye",
                        "rightIndex",
                        "rightLittle",
                        "rightRing",
                        "rightMiddle",
                        "leftIndex",
                        "leftLittle",
                        "leftRing",
                        "leftMiddle",
                        "leftThumb",
                        "rightThumb",
                        "face"
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/biometricsType"
                },
                "fullName": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{3,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "dateOfBirth": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(1869|18[7-9][0-9]|19[0-9][0-9]|20[0-9][0-9])\/([0][1-9]|1[0-2])\/([0][1-9]|[1-2][0-9]|3[01])$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "individualAuthBiometrics": {
                    "bioAttributes": [
                        "leftEye",
                        "rightEye",
                        "rightIndex",
                        "rightLittle",
                        "rightRing",
                        "rightMiddle",
                        "leftIndex",
                        "leftLittle",
                        "leftRing",
                        "leftMiddle",
                        "leftThumb",
                        "rightThumb",
                        "face"
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/biometricsType"
                },
                "introducerUIN": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "proofOfIdentity": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/documentType"
                },
                "IDSchemaVersion": {
                    "bioAttributes": [],
                    "fieldCategory": "none",
                    "format": "none",
                    "type": "number",
                    "fieldType": "default",
                    "minimum": 0
                },
                "proofOfException": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/documentType"
                },
                "phone": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^[+]*([0-9]{1})([0-9]{9})$",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "introducerName": {
                    "bioAttributes": [],
                    "fieldCategory": "evidence",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                },
                "proofOfRelationship": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/documentType"
                },
                "UIN": {
                    "bioAttributes": [],
                    "fieldCategory": "none",
                    "format": "none",
                    "type": "string",
                    "fieldType": "default"
                },
                "preferredLang": {
                    "bioAttributes": [],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "type": "string",
                    "fieldType": "dynamic"
                },
                "region": {
                    "bioAttributes": [],
                    "validators": [
                        {
                            "validator": "^(?=.{0,50}$).*",
                            "arguments": [],
                            "type": "regex"
                        }
                    ],
                    "fieldCategory": "pvt",
                    "format": "none",
                    "fieldType": "default",
                    "$ref": "#\/definitions\/simpleType"
                }
            }
        }
    }
}
'@

$ws.Range("F3").Value2 = $json
$ws.Rows.Item(3).RowHeight = 274.5

$ws.Range("H3").Select()
